# "Run `make all`" re-generation of the SSSOM schema workbook.
#
# The underlying schema did not change; the generator was simply re-run,
# which (a) re-ordered / re-named the worksheet tabs to match the LinkML
# class names ("mapping set", "mapping", "mapping registry",
# "mapping set reference"), (b) tweaked a couple of column headers on the
# "mapping set" sheet (license column renamed, mapping_tool_version column
# added), and (c) (re-)materialised the dropdown data validations that the
# generator emits for enum-typed columns on the "mapping" and
# "mapping set" sheets.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Re-order the sheets: the old "MappingSet" sheet becomes the first
#    tab; "Mapping", "MappingRegistry" and "MappingSetReference" keep
#    their relative order after it.
# ------------------------------------------------------------------
$wb.Worksheets.Item("MappingSet").Move($wb.Worksheets.Item(1))

# ------------------------------------------------------------------
# 2. Rename the tabs to their new, lower-cased/spaced names.
# ------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "mapping set"
$wb.Worksheets.Item(2).Name = "mapping"
$wb.Worksheets.Item(3).Name = "mapping registry"
$wb.Worksheets.Item(4).Name = "mapping set reference"

$wsMappingSet = $wb.Worksheets.Item("mapping set")
$wsMapping = $wb.Worksheets.Item("mapping")
$wsMappingRegistry = $wb.Worksheets.Item("mapping registry")
$wsMappingSetReference = $wb.Worksheets.Item("mapping set reference")

# ------------------------------------------------------------------
# 3. "mapping set" sheet header tweaks:
#      - column I header "mapping set_license" -> "license"
#      - new column "mapping_tool_version" inserted right after
#        "mapping_tool" (was column Q, so new column lands at R,
#        pushing mapping_date .. comment one column to the right).
# ------------------------------------------------------------------
$wsMappingSet.Range("I1").Value2 = "license"
$wsMappingSet.Range("R1").EntireColumn.Insert()
$wsMappingSet.Range("R1").Value2 = "mapping_tool_version"

# ------------------------------------------------------------------
# 4. Data validations (dropdown lists) generated for enum-typed columns.
# ------------------------------------------------------------------

$owlTypes = '"owl class,owl object property,owl data property,owl annotation property,owl named individual,skos concept,rdfs resource,rdfs class,rdfs literal,rdfs datatype,rdf property"'

# "mapping set" sheet: subject_type (J) / object_type (M)
$rng = $wsMappingSet.Range("J2:J1048576")
$rng.Validation.Add(3, 1, 1, $owlTypes)
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $false
$rng.Validation.ShowError = $false

$rng = $wsMappingSet.Range("M2:M1048576")
$rng.Validation.Add(3, 1, 1, $owlTypes)
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $false
$rng.Validation.ShowError = $false

# "mapping" sheet: predicate_modifier (F), subject_type (R), object_type (U),
# mapping_cardinality (Z)
$rng = $wsMapping.Range("F2:F1048576")
$rng.Validation.Add(3, 1, 1, '"Not"')
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $false
$rng.Validation.ShowError = $false

$rng = $wsMapping.Range("R2:R1048576")
$rng.Validation.Add(3, 1, 1, $owlTypes)
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $false
$rng.Validation.ShowError = $false

$rng = $wsMapping.Range("U2:U1048576")
$rng.Validation.Add(3, 1, 1, $owlTypes)
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $false
$rng.Validation.ShowError = $false

$rng = $wsMapping.Range("Z2:Z1048576")
$rng.Validation.Add(3, 1, 1, '"1:1,1:n,n:1,1:0,0:1,n:n"')
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $false
$rng.Validation.ShowError = $false

# ------------------------------------------------------------------
# 5. Make the first tab ("mapping set") the active one.
# ------------------------------------------------------------------
$wsMappingSet.Activate()
